# Updated cryptos list on Sun Oct  1 17:38:14 UTC 2023 with GitHub Actions
# Re-applies the scraped price/volume refresh (and the ARBITRUM / ImmutableX
# row swap at rows 36-37) cell by cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new text value. Every one of these columns (B/C/D/E) holds plain
# text in the source sheet (inline strings, no numFmt), including price
# strings that look numeric ("214.96", "0.0890", ...). Excel's normal
# Value-setter auto-coerces a numeric-looking string to a real number, which
# would both change the cell type and normalise away significant trailing/
# grouping digits (e.g. "0.0890" -> 0.089, "2.40" -> 2.4). Prefixing with a
# literal apostrophe forces text entry (same as typing '123 into a cell),
# then resetting Style back to Normal drops the implicit @ (Text) style Excel
# attaches when it text-coerces a value, so no stray style index is left on
# the cell.
$updates = [ordered]@{
    "D2" = '27.173.90'
    "D3" = '1.681.89'
    "E3" = '  +0.13%  '
    "E4" = '  +0.14%  '
    "D5" = '214.96'
    "E5" = '  -0.45%  '
    "E6" = '  -0.03%  '
    "D8" = '22.61'
    "E8" = '  +4.52%  '
    "E9" = '  +2.22%  '
    "E10" = '  +0.26%  '
    "D11" = '0.0890'
    "E11" = '  +0.06%  '
    "D12" = '1.920.97'
    "E12" = '  +0.22%  '
    "D13" = '1.683.54'
    "E13" = '  +0.33%  '
    "E14" = '  +2.15%  '
    "D15" = '0.557'
    "E15" = '  +4.76%  '
    "D16" = '66.78'
    "E16" = '  +0.48%  '
    "D17" = '27.156.72'
    "E17" = '  +0.49%  '
    "D18" = '235.95'
    "E18" = '  +0.22%  '
    "D19" = '7.88'
    "E19" = '  -3.83%  '
    "D20" = '0.0₃0739'
    "E20" = '  -0.06%  '
    "E21" = '  +0.07%  '
    "E22" = '  +1.74%  '
    "D23" = '9.53'
    "E23" = '  +2.83%  '
    "D24" = '2.08'
    "E24" = '  -1.67%  '
    "D25" = '146.88'
    "E25" = '  +0.15%  '
    "E26" = '  +2.43%  '
    "D27" = '16.34'
    "E27" = '  -0.39%  '
    "E28" = '  +0.04%  '
    "E29" = '  +0.06%  '
    "E30" = '  +1.00%  '
    "E31" = '  +0.16%  '
    "E32" = '  +0.20%  '
    "D33" = '1.544.40'
    "E33" = '  +1.71%  '
    "E34" = '  +2.25%  '
    "E35" = '  -2.82%  '
    "B36" = 'ImmutableX'
    "C36" = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    "D36" = '0.605'
    "E36" = '  +2.25%  '
    "B37" = 'ARBITRUM'
    "C37" = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    "D37" = '0.947'
    "E37" = '  +3.04%  '
    "D38" = '2.40'
    "E38" = '  -0.17%  '
    "E39" = '  -1.55%  '
    "E40" = '  +3.47%  '
    "D41" = '5.77'
    "E41" = '  +1.04%  '
    "D42" = '69.13'
    "E42" = '  +1.82%  '
    "E43" = '  +0.14%  '
    "E44" = '  -0.63%  '
    "D45" = '1.827.95'
    "E45" = '  +0.34%  '
    "E46" = '  +1.24%  '
    "D47" = '90.04'
    "E47" = '  -0.18%  '
    "E48" = '  +4.34%  '
    "D49" = '1.62'
    "E49" = '  +6.44%  '
    "D50" = '8.21'
    "E50" = '  +3.25%  '
    "E51" = '  -0.16%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = "Normal"
}
